# Auto-generated Excel COM-interop script to apply the diff changes
# described in the commit: 'Update gh-pages to output generated at 456a3b4'
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 130
$ws.Range("F6").Value = 57
$ws.Range("F8").Value = 591
$ws.Range("F9").Value = 3
$ws.Range("F11").Value = 10408
$ws.Range("F12").Value = 188
$ws.Range("F13").Value = 74
$ws.Range("F14").Value = 121
$ws.Range("F15").Value = 1987
$ws.Range("F20").Value = 57
$ws.Range("F21").Value = 226
$ws.Range("F22").Value = 1130
$ws.Range("F23").Value = 101
$ws.Range("F24").Value = 117
$ws.Range("F25").Value = 642
$ws.Range("F26").Value = 61
$ws.Range("F27").Value = 189
$ws.Range("F28").Value = 1515
$ws.Range("F29").Value = 639
$ws.Range("F30").Value = 2910
$ws.Range("F31").Value = 979
$ws.Range("F32").Value = 708
$ws.Range("F36").Value = 890
$ws.Range("F38").Value = 218
$ws.Range("B40").Value = "'2024-11-09"
$ws.Range("C40").Value = '杭州·HD·02动漫游戏博览会'
$ws.Range("E40").Value = '2024.11.09 10:00-11.10 17:00'
$ws.Range("F40").Value = 1054
$ws.Range("H40").Value = 'https://show.bilibili.com/platform/detail.html?id=92537'
$ws.Range("I40").Value = '//i1.hdslb.com/bfs/openplatform/202409/nYPayxBc1725780987986.jpeg'
$ws.Range("C41").Value = '杭州·New World动漫博览会'
$ws.Range("E41").Value = '2024.11.09 00:00-11.10 17:00'
$ws.Range("F41").Value = 5306
$ws.Range("H41").Value = 'https://show.bilibili.com/platform/detail.html?id=92136'
$ws.Range("I41").Value = '//i0.hdslb.com/bfs/openplatform/202409/63fEMuME1725960127951.jpeg'
$ws.Range("C42").Value = '杭州·岚梦国潮·夏日盛会'
$ws.Range("D42").Value = '景兴路896号 EFCLIVE欧美广场'
$ws.Range("E42").Value = '2024.11.09 10:00-11.10 18:00'
$ws.Range("F42").Value = 105
$ws.Range("G42").Value = 60
$ws.Range("H42").Value = 'https://show.bilibili.com/platform/detail.html?id=89829'
$ws.Range("I42").Value = '//i0.hdslb.com/bfs/openplatform/202407/t5Yy5W5F1721806075553.jpeg'
$ws.Range("C43").Value = '杭州·巨人only同人展中学篇'
$ws.Range("D43").Value = '康候圣街99号 顺丰创新中心'
$ws.Range("E43").Value = '2024.11.09 09:30-11.09 17:30'
$ws.Range("F43").Value = 82
$ws.Range("G43").Value = 79
$ws.Range("H43").Value = 'https://show.bilibili.com/platform/detail.html?id=92439'
$ws.Range("I43").Value = '//i2.hdslb.com/bfs/openplatform/202409/otLmkybJ1726115788486.jpeg'
$ws.Range("B44").Value = "'2024-11-10"
$ws.Range("C44").Value = '杭州·崩坏同人ONLY 爱莉希雅生日会'
$ws.Range("E44").Value = '2024.11.10 08:00-11.10 20:00'
$ws.Range("F44").Value = 119
$ws.Range("H44").Value = 'https://show.bilibili.com/platform/detail.html?id=92228'
$ws.Range("I44").Value = '//i0.hdslb.com/bfs/openplatform/202409/1FsO31h71725897488610.jpeg'
$ws.Range("B45").Value = "'2024-11-16"
$ws.Range("C45").Value = '杭州·ET金色齿轮国乙同人only'
$ws.Range("D45").Value = '转塘街道珊瑚沙东路9号 杭州白金汉爵大酒店'
$ws.Range("E45").Value = '2024.11.16 09:30-11.16 22:00'
$ws.Range("F45").Value = 200
$ws.Range("G45").Value = 25
$ws.Range("H45").Value = 'https://show.bilibili.com/platform/detail.html?id=92511'
$ws.Range("I45").Value = '//i1.hdslb.com/bfs/openplatform/202409/XfT00A611726134427042.jpeg'
$ws.Range("B46").Value = "'2024-11-23"
$ws.Range("C46").Value = '杭州·火影忍者同人only2.0 日夜连场'
$ws.Range("D46").Value = '金一路79号 XPACE湾区数字公园'
$ws.Range("E46").Value = '2024.11.23 10:00-11.23 22:30'
$ws.Range("F46").Value = 43
$ws.Range("G46").Value = 69
$ws.Range("H46").Value = 'https://show.bilibili.com/platform/detail.html?id=92097'
$ws.Range("I46").Value = '//i2.hdslb.com/bfs/openplatform/202409/q3I7lKmY1725591212982.jpeg'
$ws.Range("C47").Value = '杭州·相聚广陵代号鸢同人only3.0-三千世界'
$ws.Range("D47").Value = '康候圣街99号 顺丰创新中心'
$ws.Range("E47").Value = '2024.11.23 09:30-11.23 17:00'
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 80
$ws.Range("H47").Value = 'https://show.bilibili.com/platform/detail.html?id=92672'
$ws.Range("I47").Value = '//i0.hdslb.com/bfs/openplatform/202409/NR40ECNZ1726740199589.jpeg'
$ws.Range("F49").Value = 68
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 20
$ws.Range("F8").Value = 64
$ws.Range("F14").Value = 187
$ws = $wb.Worksheets.Item(3)
$ws.Range("F3").Value = 398
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 398
$ws.Range("F7").Value = 130
$ws.Range("F9").Value = 20
$ws.Range("F11").Value = 57
$ws.Range("F13").Value = 591
$ws.Range("F15").Value = 10408
$ws.Range("F16").Value = 74
$ws.Range("F17").Value = 121
$ws.Range("F18").Value = 1987
$ws.Range("F21").Value = 57
$ws.Range("F22").Value = 1130
$ws.Range("F23").Value = 101
$ws.Range("F24").Value = 117
$ws.Range("F26").Value = 642
$ws.Range("F27").Value = 61
$ws.Range("F28").Value = 189
$ws.Range("F29").Value = 639
$ws.Range("F30").Value = 2910
$ws.Range("F31").Value = 979
$ws.Range("F32").Value = 64
$ws.Range("F34").Value = 708
$ws.Range("F39").Value = 218
$ws.Range("F41").Value = 1062
$ws.Range("F44").Value = 119
$ws.Range("F45").Value = 200
